$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price (column D) and 1h volume-change percentage (column E)
# values, one row per coin. An empty D value means that row's price cell is
# unchanged by this update (only the volume percentage moved).
$updates = @(
    @{ Row = 2; D = "49.954.48"; E = "  +3.91%  " }
    @{ Row = 3; D = "2.651.05"; E = "  +6.24%  " }
    @{ Row = 4; D = "0.999"; E = "  +0.02%  " }
    @{ Row = 5; D = "113.98"; E = "  +7.74%  " }
    @{ Row = 6; D = "326.42"; E = "  +2.07%  " }
    @{ Row = 7; D = "0.529"; E = "  +1.69%  " }
    @{ Row = 8; D = ""; E = "  +0.03%  " }
    @{ Row = 9; D = "0.555"; E = "  +3.24%  " }
    @{ Row = 10; D = "41.02"; E = "  +5.66%  " }
    @{ Row = 11; D = "20.20"; E = "  +0.85%  " }
    @{ Row = 12; D = "0.0821"; E = "  +2.44%  " }
    @{ Row = 13; D = ""; E = "  +0.97%  " }
    @{ Row = 14; D = "7.36"; E = "  +3.86%  " }
    @{ Row = 15; D = "3.064.27"; E = "  +6.16%  " }
    @{ Row = 16; D = "2.651.33"; E = "  +6.44%  " }
    @{ Row = 17; D = "0.874"; E = "  +5.30%  " }
    @{ Row = 18; D = "49.834.85"; E = "  +4.01%  " }
    @{ Row = 19; D = "13.14"; E = "  +1.10%  " }
    @{ Row = 20; D = ""; E = "  +2.52%  " }
    @{ Row = 21; D = "2.93"; E = "  -1.61%  " }
    @{ Row = 22; D = "0.0₃0958"; E = "  +2.94%  " }
    @{ Row = 23; D = "72.08"; E = "  +1.49%  " }
    @{ Row = 24; D = "277.28"; E = "  +2.59%  " }
    @{ Row = 25; D = "2.59"; E = "  +3.00%  " }
    @{ Row = 26; D = "26.84"; E = "  +4.15%  " }
    @{ Row = 27; D = ""; E = "  -0.02%  " }
    @{ Row = 28; D = "9.99"; E = "  +2.79%  " }
    @{ Row = 29; D = ""; E = "  -2.96%  " }
    @{ Row = 30; D = "36.17"; E = "  +3.53%  " }
    @{ Row = 31; D = ""; E = "  +1.06%  " }
    @{ Row = 32; D = "50.22"; E = "  +1.84%  " }
    @{ Row = 33; D = ""; E = "  +3.25%  " }
    @{ Row = 34; D = "19.46"; E = "  +2.23%  " }
    @{ Row = 35; D = "0.0811"; E = "  +5.25%  " }
    @{ Row = 36; D = ""; E = "  -0.12%  " }
    @{ Row = 37; D = "2.08"; E = "  +7.45%  " }
    @{ Row = 38; D = ""; E = "  +5.72%  " }
    @{ Row = 39; D = ""; E = "  +8.65%  " }
    @{ Row = 40; D = "124.06"; E = "  +1.94%  " }
    @{ Row = 41; D = ""; E = "  +2.10%  " }
    @{ Row = 42; D = ""; E = "  -0.06%  " }
    @{ Row = 43; D = "22.01"; E = "  -1.45%  " }
    @{ Row = 44; D = ""; E = "  +4.26%  " }
    @{ Row = 45; D = "2.083.58"; E = "  +4.40%  " }
    @{ Row = 46; D = "3.33"; E = "  +6.33%  " }
    @{ Row = 47; D = "2.34"; E = "  +16.86%  " }
    @{ Row = 48; D = "1.98"; E = "  +5.80%  " }
    @{ Row = 49; D = "9.14"; E = "  +2.59%  " }
    @{ Row = 50; D = "5.39"; E = "  +4.78%  " }
    @{ Row = 51; D = "59.28"; E = "  +4.64%  " }
)

foreach ($u in $updates) {
    if ($u.D -ne "") {
        # Keep these price cells as plain text (matching the source data,
        # which uses "."-grouped strings like "49.954.48") rather than letting
        # Excel auto-convert numeric-looking text into a Number.
        $ws.Cells.Item($u.Row, 4).NumberFormat = "@"
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
